$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11..121 down to 12..122.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with a new weekly price observation.
# Columns A,B,C,E,F,G,H,I,J,K,L,T are identical to the (now shifted) row 12 /
# the former row 11 record; only the date + price/unit/origin columns differ.
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44515
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100112025
$ws.Range("J11").Value = "Frutilla"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 9500
$ws.Range("Q11").Value = "$/bandeja 7 kilos"
$ws.Range("R11").Value = "Provincia de Melipilla"
$ws.Range("S11").Value = 1357
$ws.Range("T11").Value = 7
